# Applies the "Trade #33 closed" update to the live trading results workbook:
#   - Summary sheet: Total Trades 32 -> 33, Win Rate % 28.12 -> 27.27
#   - Strategy Status sheet: MarketMaking row Trades 32 -> 33, Win Rate % 28.12 -> 27.27
#   - All Trades sheet: append new trade row (#33) at row 34
#   - MarketMaking sheet: append the same new trade row (#33) at row 34

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a text value into a cell without letting Excel's
# autodetection turn date-like strings (e.g. "2026-02-17") into date
# serial numbers. We temporarily force a text number format, assign the
# value, then clear the formatting again so the cell is left unstyled.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = 33
$wsSummary.Range("B9").Value = 27.27

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet (MarketMaking strategy row)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 33
$wsStatus.Range("G4").Value = 27.27

# ---------------------------------------------------------------------------
# 3) Append new trade row (#33) to both "All Trades" and "MarketMaking"
#    sheets - both currently end at row 33 (dimension A1:Q33) and gain a
#    new row 34 with identical content.
# ---------------------------------------------------------------------------
$newRow = 34

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A$newRow").Value = 33
    Set-TextValue $ws.Range("B$newRow") "2026-02-17"
    Set-TextValue $ws.Range("C$newRow") "15:22:48"
    Set-TextValue $ws.Range("D$newRow") "MarketMaking"
    Set-TextValue $ws.Range("E$newRow") "DOWN"
    $ws.Range("F$newRow").Value = 0.25
    $ws.Range("G$newRow").Value = 0.25
    Set-TextValue $ws.Range("H$newRow") "CLOSED"
    $ws.Range("I$newRow").Value = 0
    $ws.Range("J$newRow").Value = 0
    $ws.Range("K$newRow").Value = 99.76000000000001
    $ws.Range("L$newRow").Value = 0
    $ws.Range("M$newRow").Value = 0
    $ws.Range("N$newRow").Value = 0.6
    Set-TextValue $ws.Range("O$newRow") "Normal spread capture: 19600 bps"
    Set-TextValue $ws.Range("P$newRow") "early_exit"
    $ws.Range("Q$newRow").Value = 0.13
}
